$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.097.01"
$ws.Range("E2").Value = "  -0.83%  "

$ws.Range("D3").Value = "1.648.84"
$ws.Range("E3").Value = "  -1.04%  "

$ws.Range("E4").Value = "  -0.46%  "

$ws.Range("D5").Value = "'217.24"
$ws.Range("E5").Value = "  -0.96%  "

$ws.Range("D6").Value = "'0.5193"
$ws.Range("E6").Value = "  -3.02%  "

$ws.Range("E7").Value = "  -0.41%  "

$ws.Range("E8").Value = "  -1.77%  "

$ws.Range("D9").Value = "'0.06282"
$ws.Range("E9").Value = "  -2.01%  "

$ws.Range("D10").Value = "'20.40"
$ws.Range("E10").Value = "  -1.77%  "

$ws.Range("D11").Value = "'0.07783"
$ws.Range("E11").Value = "  -0.79%  "

$ws.Range("D12").Value = "1.681.57"
$ws.Range("E12").Value = "  +0.92%  "

$ws.Range("D13").Value = "'4.473"
$ws.Range("E13").Value = "  -2.09%  "

$ws.Range("D14").Value = "1.876.99"
$ws.Range("E14").Value = "  -0.91%  "

$ws.Range("D15").Value = "'0.5547"
$ws.Range("E15").Value = "  +0.11%  "

$ws.Range("D16").Value = "0.0₅7972"
$ws.Range("E16").Value = "  -3.28%  "

$ws.Range("D17").Value = "'64.71"
$ws.Range("E17").Value = "  -1.78%  "

$ws.Range("D18").Value = "26.096.97"
$ws.Range("E18").Value = "  -0.90%  "

$ws.Range("D19").Value = "'1.005"
$ws.Range("E19").Value = "  -0.48%  "

$ws.Range("D20").Value = "'4.626"
$ws.Range("E20").Value = "  -1.60%  "

$ws.Range("D21").Value = "'193.41"
$ws.Range("E21").Value = "  -0.56%  "

$ws.Range("E22").Value = "  -2.42%  "

$ws.Range("E23").Value = "  -1.79%  "

$ws.Range("E24").Value = "  -0.45%  "

$ws.Range("D25").Value = "'146.65"
$ws.Range("E25").Value = "  +0.25%  "

$ws.Range("D26").Value = "'0.1204"
$ws.Range("E26").Value = "  -2.32%  "

$ws.Range("D27").Value = "'7.158"
$ws.Range("E27").Value = "  -0.82%  "

$ws.Range("D28").Value = "'15.89"
$ws.Range("E28").Value = "  -1.88%  "

$ws.Range("D29").Value = "'1.478"
$ws.Range("E29").Value = "  -0.71%  "

$ws.Range("D30").Value = "'0.05618"
$ws.Range("E30").Value = "  -4.01%  "

$ws.Range("D31").Value = "'1.264"
$ws.Range("E31").Value = "  -1.40%  "

$ws.Range("D32").Value = "'3.480"
$ws.Range("E32").Value = "  -3.87%  "

$ws.Range("D33").Value = "'3.378"
$ws.Range("E33").Value = "  +2.88%  "

$ws.Range("D34").Value = "'1.594"
$ws.Range("E34").Value = "  -1.77%  "

$ws.Range("D35").Value = "'2.804"
$ws.Range("E35").Value = "  -0.80%  "

$ws.Range("D36").Value = "'0.9462"
$ws.Range("E36").Value = "  -2.61%  "

$ws.Range("D37").Value = "'2.404"
$ws.Range("E37").Value = "  -0.71%  "

$ws.Range("E38").Value = "  -3.09%  "

$ws.Range("D39").Value = "'5.950"
$ws.Range("E39").Value = "  +1.45%  "

$ws.Range("D40").Value = "'0.01573"
$ws.Range("E40").Value = "  -2.30%  "

$ws.Range("D41").Value = "1.061.28"
$ws.Range("E41").Value = "  +0.81%  "

$ws.Range("E42").Value = "  -0.49%  "

$ws.Range("D43").Value = "'0.8377"
$ws.Range("E43").Value = "  -4.26%  "

$ws.Range("D44").Value = "'102.93"
$ws.Range("E44").Value = "  -2.21%  "

$ws.Range("D45").Value = "1.788.48"
$ws.Range("E45").Value = "  -0.91%  "

$ws.Range("D46").Value = "'57.09"
$ws.Range("E46").Value = "  -1.59%  "

$ws.Range("D47").Value = "0.0₈105"
$ws.Range("E47").Value = "  +1.57%  "

$ws.Range("D48").Value = "'1.004"
$ws.Range("E48").Value = "  -0.91%  "

$ws.Range("D49").Value = "'0.05321"
$ws.Range("E49").Value = "  +2.97%  "

$ws.Range("D50").Value = "'0.4335"
$ws.Range("E50").Value = "  -1.13%  "

$ws.Range("D51").Value = "'7.949"
$ws.Range("E51").Value = "  -1.20%  "
